$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mattson LM"

$ws.Range("B3").Value = "Mattson LM"

$ws.Range("B4").Value = "Mattson LM"
$ws.Range("I4").Value = 78
$ws.Range("K4").Value = 95

$ws.Range("B5").Value = "Mattson LM"
$ws.Range("I5").Value = 41
$ws.Range("K5").Value = 85

$ws.Range("B6").Value = "Mattson LM"
$ws.Range("I6").Value = 28
$ws.Range("J6").Value = 34
$ws.Range("K6").Value = 29

$ws.Range("B7").Value = "Mattson LM"
$ws.Range("I7").Value = 43
$ws.Range("J7").Value = 80
$ws.Range("K7").Value = 78

$ws.Range("B8").Value = "Mattson LM"
$ws.Range("I8").Value = 18
$ws.Range("J8").Value = 67
$ws.Range("K8").Value = 52

$ws.Range("B9").Value = "Mattson LM"
$ws.Range("I9").Value = 25
$ws.Range("J9").Value = 32
$ws.Range("K9").Value = 28

$ws.Range("B10").Value = "Mattson LM"
$ws.Range("I10").Value = 15
$ws.Range("J10").Value = 19
$ws.Range("K10").Value = 18

$ws.Range("B11").Value = "Mattson LM"
$ws.Range("I11").Value = 20
$ws.Range("J11").Value = 22
$ws.Range("K11").Value = 21

$ws.Range("B12").Value = "Mattson LM"
$ws.Range("I12").Value = 25
$ws.Range("K12").Value = 68

$ws.Range("B13").Value = "Mattson LM"
$ws.Range("J13").Value = 2
$ws.Range("K13").Value = 2

$ws.Range("B14").Value = "Mattson LM"
$ws.Range("I14").Value = 2
$ws.Range("J14").Value = 3
$ws.Range("K14").Value = 3

$ws.Range("B15").Value = "Mattson LM"
$ws.Range("I15").Value = 84
$ws.Range("K15").Value = 92

$ws.Range("B16").Value = "Mattson LM"
$ws.Range("I16").Value = 96
$ws.Range("K16").Value = 97

$ws.Range("B17").Value = "Mattson LM"
$ws.Range("I17").Value = 64
$ws.Range("K17").Value = 96

$ws.Range("B18").Value = "Mattson LM"
$ws.Range("I18").Value = 2
$ws.Range("J18").Value = 3
$ws.Range("K18").Value = 3

$ws.Range("B19").Value = "Mattson LM"
$ws.Range("I19").Value = 66
$ws.Range("K19").Value = 76

$ws.Range("B20").Value = "Mattson LM"
$ws.Range("I20").Value = 9
$ws.Range("J20").Value = 20
$ws.Range("K20").Value = 18

$ws.Range("B21").Value = "Mattson LM"

$ws.Range("B22").Value = "Mattson LM"
$ws.Range("I22").Value = 63
$ws.Range("K22").Value = 80

$ws.Range("B23").Value = "Mattson LM"
$ws.Range("I23").Value = 55
$ws.Range("K23").Value = 81

$ws.Range("B24").Value = "Mattson LM"
$ws.Range("I24").Value = 63
$ws.Range("K24").Value = 85

$ws.Range("B25").Value = "Mattson LM"
$ws.Range("I25").Value = 70
$ws.Range("K25").Value = 98

$ws.Range("B26").Value = "Mattson LM"
$ws.Range("I26").Value = 10
$ws.Range("J26").Value = 75
$ws.Range("K26").Value = 55

$ws.Range("B27").Value = "Mattson LM"
$ws.Range("I27").Value = 88
$ws.Range("J27").Value = 108
$ws.Range("K27").Value = 100

$ws.Range("B28").Value = "Mattson LM"
$ws.Range("I28").Value = 119
$ws.Range("J28").Value = 136
$ws.Range("K28").Value = 124

$ws.Range("B29").Value = "Mattson LM"
$ws.Range("I29").Value = 4
$ws.Range("J29").Value = 20
$ws.Range("K29").Value = 5

$ws.Range("B30").Value = "Mattson LM"
$ws.Range("J30").Value = 2
$ws.Range("K30").Value = 2

$ws.Range("B31").Value = "Mattson LM"

$ws.Range("B32").Value = "Mattson LM"
$ws.Range("I32").Value = 69
$ws.Range("K32").Value = 74

$ws.Range("B33").Value = "Mattson LM"
$ws.Range("I33").Value = 59
$ws.Range("K33").Value = 77

$ws.Range("B34").Value = "Mattson LM"
$ws.Range("I34").Value = 55
$ws.Range("K34").Value = 96

$ws.Range("B35").Value = "Mattson LM"
$ws.Range("I35").Value = 49
$ws.Range("J35").Value = 51
$ws.Range("K35").Value = 50

$ws.Range("B36").Value = "Mattson LM"
$ws.Range("I36").Value = 121
$ws.Range("J36").Value = 132

$ws.Range("B37").Value = "Mattson LM"
$ws.Range("I37").Value = 22
$ws.Range("J37").Value = 104
$ws.Range("K37").Value = 98

$ws.Range("B38").Value = "Mattson LM"
$ws.Range("I38").Value = 60
$ws.Range("J38").Value = 135
$ws.Range("K38").Value = 122

$ws.Range("B39").Value = "Mattson LM"
$ws.Range("I39").Value = 64
$ws.Range("K39").Value = 68

$ws.Range("B40").Value = "Mattson LM"
$ws.Range("I40").Value = 32
$ws.Range("K40").Value = 53

$ws.Range("B41").Value = "Mattson LM"
$ws.Range("I41").Value = 18
$ws.Range("K41").Value = 56

$ws.Range("B42").Value = "Mattson LM"
$ws.Range("J42").Value = 53
$ws.Range("K42").Value = 50

$ws.Range("B43").Value = "Mattson LM"
$ws.Range("I43").Value = 25
$ws.Range("J43").Value = 44
$ws.Range("K43").Value = 33

$ws.Range("B44").Value = "HANWASH-Master LM"
$ws.Range("I44").Value = 33
$ws.Range("K44").Value = 73

$ws.Range("B45").Value = "HANWASH-Master LM"
$ws.Range("I45").Value = 29
$ws.Range("K45").Value = 73

$ws.Range("B46").Value = "HANWASH-Master LM"
$ws.Range("I46").Value = 18
$ws.Range("J46").Value = 124
$ws.Range("K46").Value = 29

$ws.Range("B47").Value = "HANWASH-Master LM"
$ws.Range("I47").Value = 2
$ws.Range("J47").Value = 140
$ws.Range("K47").Value = 52

$ws.Range("B48").Value = "HANWASH-Master LM"
$ws.Range("I48").Value = 12
$ws.Range("J48").Value = 193
$ws.Range("K48").Value = 170

$ws.Range("B49").Value = "HANWASH-Master LM"
$ws.Range("I49").Value = 70
$ws.Range("J49").Value = 111
$ws.Range("K49").Value = 82

$ws.Range("B50").Value = "HANWASH-Master LM"
$ws.Range("I50").Value = 98
$ws.Range("J50").Value = 151
$ws.Range("K50").Value = 102

$ws.Range("B51").Value = "HANWASH-Master LM"
$ws.Range("I51").Value = 50
$ws.Range("J51").Value = 157
$ws.Range("K51").Value = 111

$ws.Range("B52").Value = "Mattson LM"

$ws.Range("B53").Value = "Mattson LM"
$ws.Range("I53").Value = 984
$ws.Range("J53").Value = 7129
$ws.Range("K53").Value = 2605

$ws.Range("B54").Value = "Mattson LM"
$ws.Range("I54").Value = 31
$ws.Range("K54").Value = 91

$ws.Range("B55").Value = "Mattson LM"
$ws.Range("I55").Value = 6921
$ws.Range("J55").Value = 19017
$ws.Range("K55").Value = 16148

$ws.Range("B56").Value = "Mattson LM"
$ws.Range("I56").Value = 1911
$ws.Range("J56").Value = 19652
$ws.Range("K56").Value = 3023

$ws.Range("B57").Value = "Mattson LM"
$ws.Range("K57").Value = 89

$ws.Range("B58").Value = "Mattson LM"
$ws.Range("I58").Value = 9
$ws.Range("J58").Value = 30
$ws.Range("K58").Value = 11

$ws.Range("B59").Value = "Mattson LM"
$ws.Range("I59").Value = 27
$ws.Range("J59").Value = 72
$ws.Range("K59").Value = 31
